$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price (D) and 1h volume-change (E) columns.
# D-column values are leading-apostrophe prefixed so Excel keeps them as
# text (matching the source inlineStr cells) instead of auto-coercing
# plain-looking numbers like "1.00" into numeric values.

$ws.Range("D2").Value = "'64.872.82"
$ws.Range("E2").Value = '  +3.67%  '
$ws.Range("D3").Value = "'3.382.01"
$ws.Range("E3").Value = '  +3.52%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'559.37"
$ws.Range("E5").Value = '  +4.12%  '
$ws.Range("D6").Value = "'173.63"
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = '  +2.58%  '
$ws.Range("D8").Value = "'3.369.42"
$ws.Range("E8").Value = '  +3.49%  '
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = '  +12.15%  '
$ws.Range("D11").Value = "'0.629"
$ws.Range("E11").Value = '  +4.43%  '
$ws.Range("D12").Value = "'53.91"
$ws.Range("E12").Value = '  +3.50%  '
$ws.Range("D13").Value = "'0.0000278"
$ws.Range("E13").Value = '  +6.61%  '
$ws.Range("D14").Value = "'9.08"
$ws.Range("E14").Value = '  +3.99%  '
$ws.Range("D15").Value = "'3.923.19"
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("D16").Value = "'18.25"
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("D17").Value = "'0.119"
$ws.Range("E17").Value = '  +3.42%  '
$ws.Range("D18").Value = "'3.364.48"
$ws.Range("E18").Value = '  +2.28%  '
$ws.Range("D19").Value = "'64.802.96"
$ws.Range("E19").Value = '  +3.39%  '
$ws.Range("D20").Value = "'11.78"
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").Value = "'0.993"
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("D22").Value = "'469.72"
$ws.Range("E22").Value = '  +15.01%  '
$ws.Range("D23").Value = "'4.99"
$ws.Range("E23").Value = '  +14.89%  '
$ws.Range("D24").Value = "'4.13"
$ws.Range("E24").Value = '  +4.11%  '
$ws.Range("D25").Value = "'86.67"
$ws.Range("E25").Value = '  +5.79%  '
$ws.Range("D26").Value = "'13.52"
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").Value = "'2.90"
$ws.Range("E27").Value = '  +8.28%  '
$ws.Range("D28").Value = "'10.79"
$ws.Range("E28").Value = '  +3.76%  '
$ws.Range("D29").Value = "'8.74"
$ws.Range("E29").Value = '  +3.18%  '
$ws.Range("D30").Value = "'30.60"
$ws.Range("E30").Value = '  +6.97%  '
$ws.Range("D31").Value = "'6.68"
$ws.Range("E31").Value = '  +6.88%  '
$ws.Range("D32").Value = "'11.48"
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("D33").Value = "'571.27"
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").Value = "'61.35"
$ws.Range("E34").Value = '  +7.03%  '
$ws.Range("D35").Value = "'0.108"
$ws.Range("E35").Value = '  +3.10%  '
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").Value = "'3.65"
$ws.Range("E37").Value = '  +8.76%  '
$ws.Range("D38").Value = "'0.139"
$ws.Range("E38").Value = '  -3.39%  '
$ws.Range("D39").Value = "'35.56"
$ws.Range("E39").Value = '  +2.87%  '
$ws.Range("D40").Value = "'0.0₃0743"
$ws.Range("E40").Value = '  +2.59%  '
$ws.Range("D41").Value = "'0.369"
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("D42").Value = "'3.086.07"
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = "'2.84"
$ws.Range("E44").Value = '  +4.44%  '
$ws.Range("D45").Value = "'0.0415"
$ws.Range("E45").Value = '  +5.28%  '
$ws.Range("D46").Value = "'0.134"
$ws.Range("E46").Value = '  +6.16%  '
$ws.Range("D47").Value = "'2.46"
$ws.Range("E47").Value = '  +3.36%  '
$ws.Range("D48").Value = "'3.13"
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("D49").Value = "'2.59"
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("D50").Value = "'138.13"
$ws.Range("E50").Value = '  +5.18%  '
$ws.Range("D51").Value = "'8.26"
$ws.Range("E51").Value = '  +4.89%  '
